$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of MAC-address-like device data appended after the existing
# table (rows 2-146). Columns: A=regcntr_id B=machine_id C=device_id
# D=lang_code E=is_active F=cr_by G=cr_dtimes H=eff_dtimes
$rows = @(
    @(10001, 10030, 3000166),
    @(10001, 10030, 3000167),
    @(10001, 10030, 3000168),
    @(10001, 10030, 3000169),
    @(10001, 10030, 3000170),
    @(10001, 10031, 3000171),
    @(10001, 10031, 3000172),
    @(10001, 10031, 3000173),
    @(10001, 10031, 3000174),
    @(10001, 10031, 3000175)
)

$startRow = 147
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

$lastRow = $startRow + $rows.Count - 1

# Scroll the view and select the cell below the newly entered data,
# mirroring what Excel records after typing the new rows in.
$ws.Application.ActiveWindow.ScrollRow = 142
$ws.Range("A148").Select()
